$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.143982389700622
$ws.Range("B3").Value = 10.00000009998912
$ws.Range("B5").Value = 2.500000024975229
$ws.Range("B7").Value = 0.1358702519941318
$ws.Range("B8").Value = 5.974975678830996
$ws.Range("B9").Value = 0.8466759763537778
$ws.Range("B10").Value = [double]"2.627102812685198e-06"
$ws.Range("B11").Value = -0.846673349250965
$ws.Range("B12").Value = 0.03999998999529015
$ws.Range("B13").Value = 0.3999999900111381
$ws.Range("B14").Value = 0.04000000999922655
$ws.Range("B15").Value = 33027.94907862497
$ws.Range("B16").Value = 1.489473714869758
$ws.Range("B17").Value = 0.1734057146235082
$ws.Range("B18").Value = 0.04037606265197123
$ws.Range("B19").Value = 0.05165666525041572
$ws.Range("B20").Value = 0.8984967341029297
$ws.Range("B21").Value = 3.323304392746219
$ws.Range("B22").Value = -4.291134540853666
$ws.Range("B23").Value = 1.016816387783703
$ws.Range("B24").Value = -39.72418671720052
$ws.Range("B25").Value = 0.1915576723065062
$ws.Range("B26").Value = 3.946748343171597
$ws.Range("B27").Value = 0.7168557603318466
$ws.Range("B28").Value = 2.21277354021489
$ws.Range("B29").Value = 1.016816387783703
$ws.Range("B30").Value = 0.0003026548411574736
$ws.Range("B31").Value = [double]"1.506194090072323e-09"
$ws.Range("B32").Value = 0.02378068587467642
$ws.Range("B33").Value = [double]"-1.009948114489964e-17"
$ws.Range("B34").Value = 0.01992155689141272
$ws.Range("B35").Value = 0.00885673998055881
$ws.Range("B36").Value = 0.0009539459836602311
$ws.Range("B37").Value = 0.008813731172952078
$ws.Range("B38").Value = [double]"3.332735595446488e-21"
$ws.Range("B39").Value = -0.0002503571623321952
$ws.Range("B40").Value = [double]"-2.280158113993334e-20"
$ws.Range("B41").Value = 135.8702519941318
$ws.Range("B42").Value = 0.2995426311825568
$ws.Range("B43").Value = 0.04143790475840792
$ws.Range("B44").Value = 0.007336442508956208
$ws.Range("B45").Value = 0.01501426640938007
$ws.Range("B50").Value = 0.002149287324855123
$ws.Range("B52").Value = 0.02506788788185603
$ws.Range("B53").Value = 0.01006446311067643
$ws.Range("B54").Value = 0.01436147517497544
$ws.Range("B55").Value = -0.2723282682330403
$ws.Range("B56").Value = [double]"9.367102447921916e-17"
$ws.Range("B57").Value = [double]"1.250826987070088e-16"
$ws.Range("B58").Value = 1.627154779350576
$ws.Range("B59").Value = 0.06212886750254282
$ws.Range("B60").Value = 0.2723282682330403
$ws.Range("B61").Value = [double]"9.367102447921916e-17"
$ws.Range("B62").Value = [double]"2.279535067650773e-17"
$ws.Range("B63").Value = [double]"1.890877734039433e-16"
$ws.Range("B64").Value = 0.8921212093146058
$ws.Range("B65").Value = 0.03406343459952656
$ws.Range("B66").Value = 0.1493095967696309
$ws.Range("B67").Value = [double]"3.448001331619686e-17"
$ws.Range("B68").Value = [double]"1.073521418550228e-16"
$ws.Range("B69").Value = [double]"6.960262238131471e-17"
$ws.Range("B70").Value = 0.2247130993194826
$ws.Range("B71").Value = 0.04761516891355775
$ws.Range("B72").Value = 2.057193054087998
$ws.Range("B73").Value = 1.070697443330118
$ws.Range("B74").Value = 0.002730582847392099
$ws.Range("B75").Value = 0.1588345683104776
$ws.Range("B76").Value = -0.08228770158190028
$ws.Range("B77").Value = 0.05388861074945468
$ws.Range("B78").Value = 0.02843751652621496
$ws.Range("B79").Value = -0.002734123489764533
$ws.Range("B80").Value = -0.01631700644930289
$ws.Range("B81").Value = -0.7042471408726034
$ws.Range("B82").Value = -0.07999999002034769
$ws.Range("B83").Value = 0.02080043691423047
$ws.Range("B84").Value = 0.02999999000122063
$ws.Range("B85").Value = 0.087230281847203
$ws.Range("B86").Value = -0.3383973704019771
$ws.Range("B87").Value = 0.03614329957715711
$ws.Range("B88").Value = -0.1703430899706801
$ws.Range("B89").Value = -0.3412339681372011
$ws.Range("B90").Value = -0.06749694827949038
$ws.Range("B91").Value = -0.05040299875542154
$ws.Range("B92").Value = 0.4279711015071631
$ws.Range("B93").Value = 0.2477992570826815
$ws.Range("B94").Value = -0.1542028497657569
$ws.Range("B95").Value = -0.09500434585380632
$ws.Range("B96").Value = -1.144762810336675
$ws.Range("B97").Value = -0.02094373634167496
$ws.Range("B98").Value = -0.0698211701672946
$ws.Range("B99").Value = -0.001385994360586507
$ws.Range("B100").Value = 0.05934167704105264
$ws.Range("B101").Value = 0.2118940453844985
$ws.Range("B102").Value = -0.01622367756332293
$ws.Range("B103").Value = 0.003837786800247243
